$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "87.273.29"
$ws.Range("E2").Value = "  -3.11%  "
$ws.Range("D3").Value = "3.067.63"
$ws.Range("E3").Value = "  -3.74%  "
$ws.Range("E4").Value = "  +0.08%  "
$ws.Range("D5").Value = "'209.38"
$ws.Range("E5").Value = "  -2.39%  "
$ws.Range("D6").Value = "'620.87"
$ws.Range("E6").Value = "  +0.25%  "
$ws.Range("D7").Value = "'0.364"
$ws.Range("E7").Value = "  -8.15%  "
$ws.Range("D8").Value = "'0.778"
$ws.Range("E8").Value = "  +12.55%  "
$ws.Range("E9").Value = "  +0.16%  "
$ws.Range("D10").Value = "3.062.88"
$ws.Range("E10").Value = "  -3.85%  "
$ws.Range("D11").Value = "'0.579"
$ws.Range("E11").Value = "  +0.31%  "
$ws.Range("D12").Value = "'0.177"
$ws.Range("E12").Value = "  +0.10%  "
$ws.Range("D13").Value = "'0.0000236"
$ws.Range("E13").Value = "  -7.87%  "
$ws.Range("D14").Value = "'5.24"
$ws.Range("E14").Value = "  -0.56%  "
$ws.Range("D15").Value = "87.305.16"
$ws.Range("E15").Value = "  -2.82%  "
$ws.Range("D16").Value = "3.636.50"
$ws.Range("E16").Value = "  -3.51%  "
$ws.Range("D17").Value = "'31.09"
$ws.Range("E17").Value = "  -5.68%  "
$ws.Range("D18").Value = "3.069.11"
$ws.Range("E18").Value = "  -3.50%  "
$ws.Range("D19").Value = "'3.32"
$ws.Range("E19").Value = "  +0.86%  "
$ws.Range("D20").Value = "'0.0000208"
$ws.Range("E20").Value = "  +3.22%  "
$ws.Range("D21").Value = "'12.98"
$ws.Range("E21").Value = "  -3.22%  "
$ws.Range("D22").Value = "'414.86"
$ws.Range("E22").Value = "  -5.33%  "
$ws.Range("D23").Value = "'8.18"
$ws.Range("E23").Value = "  -5.18%  "
$ws.Range("D24").Value = "'4.76"
$ws.Range("E24").Value = "  -6.23%  "
$ws.Range("D25").Value = "'5.30"
$ws.Range("E25").Value = "  +2.75%  "
$ws.Range("D26").Value = "'82.04"
$ws.Range("E26").Value = "  +8.87%  "
$ws.Range("D27").Value = "'11.11"
$ws.Range("E27").Value = "  -4.82%  "
$ws.Range("D28").Value = "3.242.63"
$ws.Range("E28").Value = "  -3.10%  "
$ws.Range("D29").Value = "'1.00"
$ws.Range("E29").Value = "  +0.02%  "
$ws.Range("E30").Value = "  +0.13%  "
$ws.Range("D31").Value = "'0.149"
$ws.Range("E31").Value = "  -10.56%  "
$ws.Range("D32").Value = "'7.99"
$ws.Range("E32").Value = "  -5.42%  "
$ws.Range("B33").Value = "Bittensor"
$ws.Range("C33").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D33").Value = "'492.76"
$ws.Range("E33").Value = "  -8.16%  "
$ws.Range("B34").Value = "dogwifhat"
$ws.Range("C34").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D34").Value = "'3.62"
$ws.Range("E34").Value = "  -12.85%  "
$ws.Range("D35").Value = "'0.141"
$ws.Range("E35").Value = "  +10.56%  "
$ws.Range("D36").Value = "'6.60"
$ws.Range("E36").Value = "  -6.66%  "
$ws.Range("D37").Value = "'1.78"
$ws.Range("E37").Value = "  -4.43%  "
$ws.Range("D38").Value = "'1.23"
$ws.Range("E38").Value = "  -3.09%  "
$ws.Range("D39").Value = "'21.87"
$ws.Range("E39").Value = "  -0.95%  "
$ws.Range("E40").Value = "  -0.76%  "
$ws.Range("E41").Value = "  +0.35%  "
$ws.Range("E42").Value = "  -0.13%  "
$ws.Range("D43").Value = "'0.358"
$ws.Range("E43").Value = "  -4.48%  "
$ws.Range("B44").Value = "Monero"
$ws.Range("C44").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D44").Value = "'146.64"
$ws.Range("E44").Value = "  -2.24%  "
$ws.Range("B45").Value = "Stacks"
$ws.Range("C45").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D45").Value = "'1.80"
$ws.Range("E45").Value = "  -7.24%  "
$ws.Range("D46").Value = "'0.131"
$ws.Range("E46").Value = "  +5.26%  "
$ws.Range("D47").Value = "'43.49"
$ws.Range("E47").Value = "  -0.74%  "
$ws.Range("D48").Value = "'0.0642"
$ws.Range("E48").Value = "  +8.60%  "
$ws.Range("D49").Value = "'157.79"
$ws.Range("E49").Value = "  -8.60%  "
$ws.Range("D50").Value = "'0.701"
$ws.Range("E50").Value = "  -0.34%  "
$ws.Range("D51").Value = "'1.16"
$ws.Range("E51").Value = "  -6.14%  "
